# Generate Report for Handback
# Updates the handoff/handback timestamps for the c6e1dc52-... file (row 3)
# across the Overview, zh-cn and de-de sheets, as produced by a fresh
# handback report generation.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 3 corresponds to c6e1dc52-...md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-19 08:46:10"
$wsZhCn.Range("K3").Value = "2016-08-19 08:46:26"

# --- de-de sheet: row 3 corresponds to c6e1dc52-...md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-19 08:46:15"
$wsDeDe.Range("K3").Value = "2016-08-19 08:46:32"

# --- Overview sheet: row 3 latest HO Xliff generate date reflects the
#     newest per-locale handoff datetime for that file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-19 08:46:15"
